# Updated NATMI TPM-derived LR-pair metrics for Mst1-Mst1r (recomputed with new TPM values).
# Only the cells whose computed values change under the new TPM input are updated;
# identifier columns (A-D) and any cell whose recomputed value is unchanged are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3740666666666667
$ws.Range("H2").Value = 1.1222
$ws.Range("I2").Value = 0.2216873086880208
$ws.Range("J2").Value = 0.2216873086880208
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1457
$ws.Range("N2").Value = 0.4371
$ws.Range("O2").Value = 0.01800824077862508
$ws.Range("P2").Value = 0.01800824077862508
$ws.Range("Q2").Value = 0.05450151333333335
$ws.Range("R2").Value = 0.4905136200000001
$ws.Range("S2").Value = 0.003992198432419263
$ws.Range("T2").Value = 0.003992198432419262

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3740666666666667
$ws.Range("H3").Value = 1.1222
$ws.Range("I3").Value = 0.2216873086880208
$ws.Range("J3").Value = 0.2216873086880208
$ws.Range("O3").Value = 0.6993369086413642
$ws.Range("P3").Value = 0.6993369086413641
$ws.Range("Q3").Value = 2.116526556888889
$ws.Range("R3").Value = 19.048739012
$ws.Range("S3").Value = 0.1550341171429043
$ws.Range("T3").Value = 0.1550341171429043

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3740666666666667
$ws.Range("H4").Value = 1.1222
$ws.Range("I4").Value = 0.2216873086880208
$ws.Range("J4").Value = 0.2216873086880208
$ws.Range("O4").Value = 0.2826548505800108
$ws.Range("P4").Value = 0.2826548505800108
$ws.Range("Q4").Value = 0.8554481971333335
$ws.Range("R4").Value = 7.699033774200001
$ws.Range("S4").Value = 0.06266099311269727
$ws.Range("T4").Value = 0.06266099311269725

# Row 5
$ws.Range("I5").Value = 0.6019302323054629
$ws.Range("J5").Value = 0.6019302323054631
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1457
$ws.Range("N5").Value = 0.4371
$ws.Range("O5").Value = 0.01800824077862508
$ws.Range("P5").Value = 0.01800824077862508
$ws.Range("Q5").Value = 0.1479837018
$ws.Range("R5").Value = 1.3318533162
$ws.Range("S5").Value = 0.01083970455529051
$ws.Range("T5").Value = 0.01083970455529051

# Row 6
$ws.Range("I6").Value = 0.6019302323054629
$ws.Range("J6").Value = 0.6019302323054631
$ws.Range("O6").Value = 0.6993369086413642
$ws.Range("P6").Value = 0.6993369086413641
$ws.Range("S6").Value = 0.4209520278782807
$ws.Range("T6").Value = 0.4209520278782807

# Row 7
$ws.Range("I7").Value = 0.6019302323054629
$ws.Range("J7").Value = 0.6019302323054631
$ws.Range("O7").Value = 0.2826548505800108
$ws.Range("P7").Value = 0.2826548505800108
$ws.Range("S7").Value = 0.1701384998718918
$ws.Range("T7").Value = 0.1701384998718918

# Row 8
$ws.Range("I8").Value = 0.1763824590065161
$ws.Range("J8").Value = 0.1763824590065161
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1457
$ws.Range("N8").Value = 0.4371
$ws.Range("O8").Value = 0.01800824077862508
$ws.Range("P8").Value = 0.01800824077862508
$ws.Range("Q8").Value = 0.04336337970000001
$ws.Range("R8").Value = 0.3902704173000001
$ws.Range("S8").Value = 0.00317633779091531
$ws.Range("T8").Value = 0.003176337790915309

# Row 9
$ws.Range("I9").Value = 0.1763824590065161
$ws.Range("J9").Value = 0.1763824590065161
$ws.Range("O9").Value = 0.6993369086413642
$ws.Range("P9").Value = 0.6993369086413641
$ws.Range("S9").Value = 0.1233507636201791
$ws.Range("T9").Value = 0.1233507636201791

# Row 10
$ws.Range("I10").Value = 0.1763824590065161
$ws.Range("J10").Value = 0.1763824590065161
$ws.Range("O10").Value = 0.2826548505800108
$ws.Range("P10").Value = 0.2826548505800108
$ws.Range("S10").Value = 0.04985535759542169
$ws.Range("T10").Value = 0.04985535759542168
